$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date text values that replace the former date-serial values in columns B and C
# (time / date-ness removed from the table -- the values become plain text strings)
$startDates = @(
  "2021-01-25","2021-02-25","2021-03-23","2021-04-25","2021-05-25",
  "2021-06-21","2021-07-23","2021-08-23","2021-09-25","2021-10-27",
  "2021-11-28","2021-12-27","2022-01-31","2022-03-05","2022-04-02",
  "2022-05-01","2022-06-04","2022-07-02","2022-07-30","2022-08-28",
  "2022-09-28","2022-10-31","2022-12-01","2023-01-05","2023-02-01"
)
$endDates = @(
  "2021-01-30","2021-03-03","2021-03-29","2021-05-01","2021-05-31",
  "2021-06-27","2021-07-29","2021-08-29","2021-09-30","2021-11-02",
  "2021-12-03","2022-01-01","2022-02-06","2022-03-10","2022-04-07",
  "2022-05-07","2022-06-09","2022-07-08","2022-08-04","2022-09-02",
  "2022-10-04","2022-11-06","2022-12-07","2023-01-10","2023-02-06"
)

for ($i = 0; $i -lt 25; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 2).Value = "'" + $startDates[$i]
  $ws.Cells.Item($row, 3).Value = "'" + $endDates[$i]
}

# New column for future data -- width added, no values
$ws.Columns("D").ColumnWidth = 15.7265625

# View state: zoomed out a bit, topLeftCell reset, selection moved to E26
$excel.ActiveWindow.Zoom = 73
$ws.Range("E26").Select() | Out-Null
